$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'white athletic leggings mens'
$ws.Cells.Item(2, 1).Value = 'supreme basketball tights'
$ws.Cells.Item(3, 1).Value = 'snowmobiling knee pads'
$ws.Cells.Item(4, 1).Value = 'photographer knee pads'
$ws.Cells.Item(5, 1).Value = 'raymens leggings'
$ws.Cells.Item(6, 1).Value = 'training basketball youth'
$ws.Cells.Item(7, 1).Value = 'baleaf men compression pants'
$ws.Cells.Item(8, 1).Value = 'winter leggings men'
$ws.Cells.Item(9, 1).Value = 'kids pants with knee pads'
$ws.Cells.Item(10, 1).Value = 'eclipse knee pads'
$ws.Cells.Item(11, 1).Value = 'pilates knee pads'
$ws.Cells.Item(12, 1).Value = 'fitted mens tights'
$ws.Cells.Item(13, 1).Value = 'mens knee pads bmx'
$ws.Cells.Item(14, 1).Value = 'canoe knee pad'
$ws.Cells.Item(15, 1).Value = 'graduated compression pants'
$ws.Cells.Item(16, 1).Value = 'graduated compression leggings'
$ws.Cells.Item(17, 1).Value = 'soft knee pads'
$ws.Cells.Item(18, 1).Value = 'knee pads nike'
$ws.Cells.Item(19, 1).Value = 'knee pad toddler'
$ws.Cells.Item(20, 1).Value = 'knee pad yoga'
$ws.Cells.Item(21, 1).Value = 'knee pads dance'
$ws.Cells.Item(22, 1).Value = 'knee pads mizuno'
$ws.Cells.Item(23, 1).Value = 'elbow knee pads'
$ws.Cells.Item(24, 1).Value = 'nee pads basketball'
$ws.Cells.Item(25, 1).Value = 'compression pants set men'
$ws.Cells.Item(26, 1).Value = 'mens basketball snap pants'
$ws.Cells.Item(27, 1).Value = 'nike compression pants for youth'
$ws.Cells.Item(28, 1).Value = 'mens compression tights cold weather'
$ws.Cells.Item(29, 1).Value = 'mens under armour basketball tights'
$ws.Cells.Item(30, 1).Value = 'mens tights pockets'
$ws.Cells.Item(31, 1).Value = 'rollerblades knee pads'
$ws.Cells.Item(32, 1).Value = 'knee pads 3xl'
$ws.Cells.Item(33, 1).Value = 'knee pads 8'
$ws.Cells.Item(34, 1).Value = 'knee pad hard'
$ws.Cells.Item(35, 1).Value = 'men leggings fleece'
$ws.Cells.Item(36, 1).Value = 'teflex knee pads'
$ws.Cells.Item(37, 1).Value = 'man leggings thermal'
$ws.Cells.Item(38, 1).Value = 'sailing knee pad'
$ws.Cells.Item(39, 1).Value = 'knee pads downhill'
$ws.Cells.Item(40, 1).Value = 'knee pads airsoft'
$ws.Cells.Item(41, 1).Value = 'knee pads army'
$ws.Cells.Item(42, 1).Value = 'knee pads enduro'
$ws.Cells.Item(43, 1).Value = 'knee pads bike'
$ws.Cells.Item(44, 1).Value = 'knee pads caving'
$ws.Cells.Item(45, 1).Value = 'knee pads for women'
$ws.Cells.Item(46, 1).Value = 'knee pads gloves'
$ws.Cells.Item(47, 1).Value = 'knee pads green'
$ws.Cells.Item(48, 1).Value = 'knee pads kali'
$ws.Cells.Item(49, 1).Value = 'knee pads kuangmi'
$ws.Cells.Item(50, 1).Value = 'knee pads longboard'
$ws.Cells.Item(51, 1).Value = 'knee pads neoprene'
$ws.Cells.Item(52, 1).Value = 'knee pads ocp'
$ws.Cells.Item(53, 1).Value = 'knee pads orange'
$ws.Cells.Item(54, 1).Value = 'knee pads over pants'
$ws.Cells.Item(55, 1).Value = 'knee pads plastic'
$ws.Cells.Item(56, 1).Value = 'knee pads purple'
$ws.Cells.Item(57, 1).Value = 'knee pads razor'
$ws.Cells.Item(58, 1).Value = 'knee pads red'
$ws.Cells.Item(59, 1).Value = 'knee pads rollerblading'
$ws.Cells.Item(60, 1).Value = 'knee pads sailing'
$ws.Cells.Item(61, 1).Value = 'knee pads scooter'
$ws.Cells.Item(62, 1).Value = 'knee pads set'
$ws.Cells.Item(63, 1).Value = 'knee pads shooting'
$ws.Cells.Item(64, 1).Value = 'knee pads swat'
$ws.Cells.Item(65, 1).Value = 'knee pads teen'
$ws.Cells.Item(66, 1).Value = 'knee pads tsg'
$ws.Cells.Item(67, 1).Value = 'knee pads usmc'
$ws.Cells.Item(68, 1).Value = 'knee pads viper'
$ws.Cells.Item(69, 1).Value = 'knee pads white'
$ws.Cells.Item(70, 1).Value = 'knee pads yellow'
$ws.Cells.Item(71, 1).Value = 'mens leggings xs'
$ws.Cells.Item(72, 1).Value = 'xtextile compression pants men'
$ws.Cells.Item(73, 1).Value = 'crx men''s tights'
$ws.Cells.Item(74, 1).Value = 'yoga capri pants'
$ws.Cells.Item(75, 1).Value = 'dodoing kneepads'
$ws.Cells.Item(76, 1).Value = 'cavaliers basketball leggings'
$ws.Cells.Item(77, 1).Value = 'basketball knee pads kids'
$ws.Cells.Item(78, 1).Value = 'knee pads for toddlers'
$ws.Cells.Item(79, 1).Value = 'knee pads skating'
$ws.Cells.Item(80, 1).Value = 'knee pads skateboarding'
$ws.Cells.Item(81, 1).Value = 'basketball knee pads kids boys'
$ws.Cells.Item(82, 1).Value = 'knee pads for dance'
$ws.Cells.Item(83, 1).Value = 'knee pads rollerblade'
$ws.Cells.Item(84, 1).Value = 'knee pads tan'
$ws.Cells.Item(85, 1).Value = 'knee pad dancer'
$ws.Cells.Item(86, 1).Value = 'knee pads adidas'
$ws.Cells.Item(87, 1).Value = 'knee pads basketball mcdavid'
$ws.Cells.Item(88, 1).Value = 'knee pads dancing'
$ws.Cells.Item(89, 1).Value = 'knee pads dodgeball'
$ws.Cells.Item(90, 1).Value = 'knee pads pair'
$ws.Cells.Item(91, 1).Value = 'knee pads longboarding'
$ws.Cells.Item(92, 1).Value = 'knee pads nba'
$ws.Cells.Item(93, 1).Value = 'knee pads pole'
$ws.Cells.Item(94, 1).Value = 'knee pad and elbow pads'
$ws.Cells.Item(95, 1).Value = 'knee pad adidas'
$ws.Cells.Item(96, 1).Value = 'knee pad asics'
$ws.Cells.Item(97, 1).Value = 'knee pad for kids'
$ws.Cells.Item(98, 1).Value = 'knee pad military'
$ws.Cells.Item(99, 1).Value = 'knee pad mma'
$ws.Cells.Item(100, 1).Value = 'knee pad mizuno'
